$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header labels for columns T, U, V (row 1): "Ea", "Eb", "delta Eab"
$ws.Range("T1").Value = "Ea"
$ws.Range("U1").Value = "Eb"
$ws.Range("V1").Value = "ΔEab"

# Row 3 is already hidden; temporarily unhide it so that writing new cell
# values into it does not trigger an unwanted auto row-height recalculation.
$ws.Rows.Item(3).Hidden = $false

# Fill the placeholder "todo" columns D:G and R:V for data rows 2-7
$ws.Range("D2:G7").Value = "todo"
$ws.Range("R2:V7").Value = "todo"

# Those placeholder cells should wrap text (matches the formatting used
# elsewhere in the sheet for this kind of placeholder cell)
$ws.Range("D2:G7").WrapText = $true
$ws.Range("R2:V7").WrapText = $true

# Row height adjustments
$ws.Rows.Item(2).RowHeight = 12
$ws.Rows.Item(4).RowHeight = 27
$ws.Rows.Item(6).RowHeight = 14.25

# Final hidden state: rows 3, 4 and 5 end up hidden
$ws.Rows.Item(3).Hidden = $true
$ws.Rows.Item(4).Hidden = $true
$ws.Rows.Item(5).Hidden = $true
